$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "dbExcel" column (B), shifting
# "dbExcel" -> C and "WebExcel" -> D (and their row-2 values along with
# them). The new column B holds the "StatQuery" header/value pair used
# to validate the trials stat bar.
$ws.Columns("B").Insert()

$ws.Range("B1").Value = "StatQuery"
$ws.Range("B2").Value = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN ['Cholangiocarcinoma, intrahepatic and extrahepatic bile ducts (adenocarcinoma)'] OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"

# Match the wrapped-text look of the "query" column (A2) for the new
# "StatQuery" cell (B2).
$ws.Range("B2").WrapText = $true

# New column should be the same width as column A.
$ws.Columns("B").ColumnWidth = $ws.Columns("A").ColumnWidth

$ws.Range("A3").Select()
